$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.950.18'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').Value = '3.042.12'
$ws.Range('E3').Value = '  -1.55%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '''580.81'
$ws.Range('E5').Value = '  -1.87%  '
$ws.Range('D6').Value = '''150.74'
$ws.Range('E6').Value = '  -2.85%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('D8').Value = '''0.529'
$ws.Range('E8').Value = '  -2.89%  '
$ws.Range('D9').Value = '3.044.38'
$ws.Range('E9').Value = '  -1.28%  '
$ws.Range('E10').Value = '  -3.32%  '
$ws.Range('D11').Value = '''5.77'
$ws.Range('E11').Value = '  -1.35%  '
$ws.Range('E12').Value = '  -2.33%  '
$ws.Range('E13').Value = '  -3.67%  '
$ws.Range('D14').Value = '''35.83'
$ws.Range('E14').Value = '  -4.57%  '
$ws.Range('E15').Value = '  +2.04%  '
$ws.Range('D16').Value = '3.549.01'
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('D17').Value = '''7.08'
$ws.Range('E17').Value = '  -1.56%  '
$ws.Range('D18').Value = '62.953.74'
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('D19').Value = '3.047.77'
$ws.Range('D20').Value = '''479.23'
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('E21').Value = '  -3.08%  '
$ws.Range('D22').Value = '''0.700'
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('D23').Value = '''7.47'
$ws.Range('E23').Value = '  -1.09%  '
$ws.Range('D24').Value = '''2.36'
$ws.Range('E24').Value = '  -2.21%  '
$ws.Range('D25').Value = '''81.34'
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('D26').Value = '''12.56'
$ws.Range('E26').Value = '  -2.94%  '
$ws.Range('D27').Value = '''10.52'
$ws.Range('E27').Value = '  +4.49%  '
$ws.Range('D28').Value = '''1.00'
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('D29').Value = '''7.30'
$ws.Range('E29').Value = '  -1.15%  '
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').Value = '''2.63'
$ws.Range('D32').Value = '''2.19'
$ws.Range('E32').Value = '  -0.02%  '
$ws.Range('D33').Value = '''27.56'
$ws.Range('E33').Value = '  +1.00%  '
$ws.Range('D34').Value = '''0.109'
$ws.Range('E34').Value = '  -3.88%  '
$ws.Range('E35').Value = '  +0.36%  '
$ws.Range('E36').Value = '  -5.72%  '
$ws.Range('D37').Value = '''5.84'
$ws.Range('E37').Value = '  -4.20%  '
$ws.Range('E38').Value = '  -2.64%  '
$ws.Range('D39').Value = '''3.09'
$ws.Range('E39').Value = '  -8.85%  '
$ws.Range('D40').Value = '''50.24'
$ws.Range('E40').Value = '  -0.86%  '
$ws.Range('E41').Value = '  -2.25%  '
$ws.Range('D42').Value = '''423.85'
$ws.Range('E42').Value = '  -4.42%  '
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').Value = '''0.284'
$ws.Range('E43').Value = '  -0.21%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').Value = '''0.114'
$ws.Range('E44').Value = '  +2.98%  '
$ws.Range('D45').Value = '2.829.13'
$ws.Range('E45').Value = '  +0.90%  '
$ws.Range('D46').Value = '''0.0359'
$ws.Range('E46').Value = '  -1.12%  '
$ws.Range('D47').Value = '''37.76'
$ws.Range('E47').Value = '  -5.76%  '
$ws.Range('D48').Value = '''126.77'
$ws.Range('D50').Value = '''24.79'
$ws.Range('E50').Value = '  -3.92%  '
$ws.Range('E51').Value = '  -1.46%  '
